$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting everything down
$ws.Rows.Item(1).Insert()

# Fill in the new header row
$ws.Range("A1").Value = "Hardware"
$ws.Range("B1").Value = "Mac : 2.4 i7 4 GB DDR"
$ws.Range("C1").Value = "Web"

# Move the selection/active cell to D1 as in the diff
$ws.Range("D1").Select()
